$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right 5 -> 4, Wrong -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total"): Right 105 -> 84, Wrong -5 -> -10, and the score text
$ws.Range("B12").Value = 84
$ws.Range("C12").Value = -10
$ws.Range("E12").Value = "74 / 112"
